$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data rows (60-65) appended to the dataset, matching the
# existing table layout: Mercado ID, Mercado, Region, Fecha, Codreg,
# Categoria ID, Categoria, Variedad, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg, Kg o Unidades, Clasificacion

$rows = @(
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 300000000, "Espárragos", "Sin especificar", "Banquete", 990,  1400, 1500, 1443, "$/kilo", "Provincia de Linares",  1443, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 300000000, "Espárragos", "Sin especificar", "Banquete", 660,  1300, 1400, 1335, "$/kilo", "Región Metropolitana",  1335, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 300000000, "Espárragos", "Sin especificar", "Primera",  1000, 1200, 1300, 1253, "$/kilo", "Provincia de Linares",  1253, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 300000000, "Espárragos", "Sin especificar", "Primera",  640,  1100, 1200, 1142, "$/kilo", "Región Metropolitana",  1142, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 300000000, "Espárragos", "Sin especificar", "Segunda",  630,  1000, 1000, 1000, "$/kilo", "Provincia de Linares",  1000, 1, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 300000000, "Espárragos", "Sin especificar", "Segunda",  440,  900,  1000, 941,  "$/kilo", "Región Metropolitana",  941,  1, "Hortaliza")
)

$startRow = 60
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.Value = $data[$c]
    }
    # Column D (4) holds the date serial; reuse the same date style used
    # by the rest of the column (numFmtId 165 "YYYY-MM-DD HH:MM:SS")
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
